# Daily attendance processing - 2026-01-12 09:51:13
# Applies the session-analysis updates to the "Session Analysis Results" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Summary statistics block (K5:L10)
# ---------------------------------------------------------------------------
$ws.Range("L6").Value = 237   # Recorded Sessions
$ws.Range("L7").Value = 27    # Missing Sessions

$ws.Range("L9").NumberFormat = "@"
$ws.Range("L9").Value = "74.5%"   # Coverage %

$ws.Range("L10").NumberFormat = "@"
$ws.Range("L10").Value = "75.8%"  # Average Attendance %

# ---------------------------------------------------------------------------
# 2. "Recorded By" cells: reorder "System, dnasr281@gmail.com"
#    -> "dnasr281@gmail.com, System"
# ---------------------------------------------------------------------------
$ws.Range("G8").Value = "dnasr281@gmail.com, System"
$ws.Range("G9").Value = "dnasr281@gmail.com, System"
$ws.Range("G10").Value = "dnasr281@gmail.com, System"
$ws.Range("G12").Value = "dnasr281@gmail.com, System"
$ws.Range("G14").Value = "dnasr281@gmail.com, System"
$ws.Range("G15").Value = "dnasr281@gmail.com, System"
$ws.Range("G17").Value = "dnasr281@gmail.com, System"
$ws.Range("G18").Value = "dnasr281@gmail.com, System"
$ws.Range("G34").Value = "dnasr281@gmail.com, System"
$ws.Range("G35").Value = "dnasr281@gmail.com, System"
$ws.Range("G36").Value = "dnasr281@gmail.com, System"
$ws.Range("G38").Value = "dnasr281@gmail.com, System"
$ws.Range("G40").Value = "dnasr281@gmail.com, System"
$ws.Range("G41").Value = "dnasr281@gmail.com, System"
$ws.Range("G43").Value = "dnasr281@gmail.com, System"
$ws.Range("G44").Value = "dnasr281@gmail.com, System"
$ws.Range("G60").Value = "dnasr281@gmail.com, System"
$ws.Range("G61").Value = "dnasr281@gmail.com, System"
$ws.Range("G62").Value = "dnasr281@gmail.com, System"
$ws.Range("G64").Value = "dnasr281@gmail.com, System"
$ws.Range("G66").Value = "dnasr281@gmail.com, System"
$ws.Range("G67").Value = "dnasr281@gmail.com, System"
$ws.Range("G69").Value = "dnasr281@gmail.com, System"
$ws.Range("G70").Value = "dnasr281@gmail.com, System"
$ws.Range("G86").Value = "dnasr281@gmail.com, System"
$ws.Range("G87").Value = "dnasr281@gmail.com, System"
$ws.Range("G88").Value = "dnasr281@gmail.com, System"
$ws.Range("G90").Value = "dnasr281@gmail.com, System"
$ws.Range("G92").Value = "dnasr281@gmail.com, System"
$ws.Range("G93").Value = "dnasr281@gmail.com, System"
$ws.Range("G95").Value = "dnasr281@gmail.com, System"
$ws.Range("G96").Value = "dnasr281@gmail.com, System"
$ws.Range("G112").Value = "dnasr281@gmail.com, System"
$ws.Range("G113").Value = "dnasr281@gmail.com, System"
$ws.Range("G114").Value = "dnasr281@gmail.com, System"
$ws.Range("G116").Value = "dnasr281@gmail.com, System"
$ws.Range("G118").Value = "dnasr281@gmail.com, System"
$ws.Range("G119").Value = "dnasr281@gmail.com, System"
$ws.Range("G121").Value = "dnasr281@gmail.com, System"
$ws.Range("G122").Value = "dnasr281@gmail.com, System"
$ws.Range("G138").Value = "dnasr281@gmail.com, System"
$ws.Range("G139").Value = "dnasr281@gmail.com, System"
$ws.Range("G140").Value = "dnasr281@gmail.com, System"
$ws.Range("G142").Value = "dnasr281@gmail.com, System"
$ws.Range("G144").Value = "dnasr281@gmail.com, System"
$ws.Range("G145").Value = "dnasr281@gmail.com, System"
$ws.Range("G147").Value = "dnasr281@gmail.com, System"
$ws.Range("G148").Value = "dnasr281@gmail.com, System"
$ws.Range("G164").Value = "dnasr281@gmail.com, System"
$ws.Range("G167").Value = "dnasr281@gmail.com, System"
$ws.Range("G170").Value = "dnasr281@gmail.com, System"
$ws.Range("G174").Value = "dnasr281@gmail.com, System"
$ws.Range("G191").Value = "dnasr281@gmail.com, System"
$ws.Range("G194").Value = "dnasr281@gmail.com, System"
$ws.Range("G197").Value = "dnasr281@gmail.com, System"
$ws.Range("G201").Value = "dnasr281@gmail.com, System"
$ws.Range("G218").Value = "dnasr281@gmail.com, System"
$ws.Range("G221").Value = "dnasr281@gmail.com, System"
$ws.Range("G224").Value = "dnasr281@gmail.com, System"
$ws.Range("G228").Value = "dnasr281@gmail.com, System"
$ws.Range("G245").Value = "dnasr281@gmail.com, System"
$ws.Range("G248").Value = "dnasr281@gmail.com, System"
$ws.Range("G251").Value = "dnasr281@gmail.com, System"
$ws.Range("G255").Value = "dnasr281@gmail.com, System"
$ws.Range("G272").Value = "dnasr281@gmail.com, System"
$ws.Range("G275").Value = "dnasr281@gmail.com, System"
$ws.Range("G278").Value = "dnasr281@gmail.com, System"
$ws.Range("G282").Value = "dnasr281@gmail.com, System"
$ws.Range("G299").Value = "dnasr281@gmail.com, System"
$ws.Range("G302").Value = "dnasr281@gmail.com, System"
$ws.Range("G305").Value = "dnasr281@gmail.com, System"
$ws.Range("G309").Value = "dnasr281@gmail.com, System"

# ---------------------------------------------------------------------------
# 3. Weekly per-group breakdown (rows 15-20): Recorded/Missing counts + %
# ---------------------------------------------------------------------------
$ws.Range("O15").Value = 21
$ws.Range("P15").Value = 2
$ws.Range("R15").NumberFormat = "@"
$ws.Range("R15").Value = "80.8%"
$ws.Range("S15").NumberFormat = "@"
$ws.Range("S15").Value = "79.9%"

$ws.Range("O16").Value = 22
$ws.Range("P16").Value = 1
$ws.Range("R16").NumberFormat = "@"
$ws.Range("R16").Value = "84.6%"
$ws.Range("S16").NumberFormat = "@"
$ws.Range("S16").Value = "80.3%"

$ws.Range("O17").Value = 22
$ws.Range("P17").Value = 1
$ws.Range("R17").NumberFormat = "@"
$ws.Range("R17").Value = "84.6%"
$ws.Range("S17").NumberFormat = "@"
$ws.Range("S17").Value = "72.6%"

$ws.Range("O18").Value = 22
$ws.Range("P18").Value = 1
$ws.Range("R18").NumberFormat = "@"
$ws.Range("R18").Value = "84.6%"
$ws.Range("S18").NumberFormat = "@"
$ws.Range("S18").Value = "77.3%"

$ws.Range("O19").Value = 22
$ws.Range("P19").Value = 1
$ws.Range("R19").NumberFormat = "@"
$ws.Range("R19").Value = "84.6%"
$ws.Range("S19").NumberFormat = "@"
$ws.Range("S19").Value = "77.0%"

$ws.Range("O20").Value = 21
$ws.Range("P20").Value = 2
$ws.Range("R20").NumberFormat = "@"
$ws.Range("R20").Value = "80.8%"
$ws.Range("S20").NumberFormat = "@"
$ws.Range("S20").Value = "79.1%"

# ---------------------------------------------------------------------------
# 4. Session-24 rows (24, 50, 76, 102, 128, 154) flip from "Not Recorded"
#    (pink highlight, style index 4) to "Recorded" (green highlight, style
#    index 2) now that attendance has been taken.
# ---------------------------------------------------------------------------

# Row 24 - B1A1
$ws.Range("A15:I15").Copy() | Out-Null
$ws.Range("A24:I24").PasteSpecial(-4122) | Out-Null
$ws.Range("G24").Value = "dnasr281@gmail.com"
$ws.Range("H24").Value = "16/26"
$ws.Range("I24").Value = "Recorded"

# Row 50 - B1A2
$ws.Range("A15:I15").Copy() | Out-Null
$ws.Range("A50:I50").PasteSpecial(-4122) | Out-Null
$ws.Range("G50").Value = "dnasr281@gmail.com"
$ws.Range("H50").Value = "17/27"
$ws.Range("I50").Value = "Recorded"

# Row 76 - B1B1
$ws.Range("A15:I15").Copy() | Out-Null
$ws.Range("A76:I76").PasteSpecial(-4122) | Out-Null
$ws.Range("G76").Value = "dnasr281@gmail.com"
$ws.Range("H76").Value = "21/26"
$ws.Range("I76").Value = "Recorded"

# Row 102 - B1B2
$ws.Range("A15:I15").Copy() | Out-Null
$ws.Range("A102:I102").PasteSpecial(-4122) | Out-Null
$ws.Range("G102").Value = "dnasr281@gmail.com"
$ws.Range("H102").Value = "16/27"
$ws.Range("I102").Value = "Recorded"

# Row 128 - B1C1
$ws.Range("A15:I15").Copy() | Out-Null
$ws.Range("A128:I128").PasteSpecial(-4122) | Out-Null
$ws.Range("G128").Value = "dnasr281@gmail.com"
$ws.Range("H128").Value = "21/30"
$ws.Range("I128").Value = "Recorded"

# Row 154 - B1C2
$ws.Range("A15:I15").Copy() | Out-Null
$ws.Range("A154:I154").PasteSpecial(-4122) | Out-Null
$ws.Range("G154").Value = "dnasr281@gmail.com"
$ws.Range("H154").Value = "17/23"
$ws.Range("I154").Value = "Recorded"

$excel.CutCopyMode = 0
